$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resolve earlier detect_structure flags now that later data confirms the pattern
$ws.Cells.Item(54, 17).Value = 0
$ws.Cells.Item(58, 17).Value = 0
$ws.Cells.Item(65, 17).Value = 0
$ws.Cells.Item(74, 17).Value = 0

# Flag row 167 as a pivot and backfill the backup column for the last two existing rows
$ws.Cells.Item(167, 15).Value = 1
$ws.Cells.Item(167, 18).Value = 0
$ws.Cells.Item(168, 18).Value = 0

# Append newly completed weekly bars (rows 169-191)
# Row 169
$ws.Cells.Item(169, 1).Value = 45474
$ws.Cells.Item(169, 2).Value = 1510
$ws.Cells.Item(169, 3).Value = 1538.75
$ws.Cells.Item(169, 4).Value = 1468
$ws.Cells.Item(169, 5).Value = 1533.900024414062
$ws.Cells.Item(169, 6).Value = 1531.212280273438
$ws.Cells.Item(169, 7).Value = 2108468
$ws.Cells.Item(169, 8).Value = 2024
$ws.Cells.Item(169, 9).Value = 7
$ws.Cells.Item(169, 10).Value = 1
$ws.Cells.Item(169, 11).Value = 0
$ws.Cells.Item(169, 12).Value = 0
$ws.Cells.Item(169, 13).Value = 0
$ws.Cells.Item(169, 14).Value = 27
$ws.Cells.Item(169, 15).Value = 0
$ws.Cells.Item(169, 16).Value = 0
$ws.Cells.Item(169, 17).Value = 0
$ws.Cells.Item(169, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 170
$ws.Cells.Item(170, 1).Value = 45481
$ws.Cells.Item(170, 2).Value = 1548
$ws.Cells.Item(170, 3).Value = 1595
$ws.Cells.Item(170, 4).Value = 1387.900024414062
$ws.Cells.Item(170, 5).Value = 1395.800048828125
$ws.Cells.Item(170, 6).Value = 1393.354248046875
$ws.Cells.Item(170, 7).Value = 4749170
$ws.Cells.Item(170, 8).Value = 2024
$ws.Cells.Item(170, 9).Value = 7
$ws.Cells.Item(170, 10).Value = 8
$ws.Cells.Item(170, 11).Value = 0
$ws.Cells.Item(170, 12).Value = 0
$ws.Cells.Item(170, 13).Value = 0
$ws.Cells.Item(170, 14).Value = 28
$ws.Cells.Item(170, 15).Value = 0
$ws.Cells.Item(170, 16).Value = 0
$ws.Cells.Item(170, 17).Value = 0
$ws.Cells.Item(170, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 171
$ws.Cells.Item(171, 1).Value = 45488
$ws.Cells.Item(171, 2).Value = 1416.400024414062
$ws.Cells.Item(171, 3).Value = 1522.75
$ws.Cells.Item(171, 4).Value = 1370
$ws.Cells.Item(171, 5).Value = 1448.650024414062
$ws.Cells.Item(171, 6).Value = 1446.111572265625
$ws.Cells.Item(171, 7).Value = 5862068
$ws.Cells.Item(171, 8).Value = 2024
$ws.Cells.Item(171, 9).Value = 7
$ws.Cells.Item(171, 10).Value = 15
$ws.Cells.Item(171, 11).Value = 0
$ws.Cells.Item(171, 12).Value = 0
$ws.Cells.Item(171, 13).Value = 0
$ws.Cells.Item(171, 14).Value = 29
$ws.Cells.Item(171, 15).Value = 0
$ws.Cells.Item(171, 16).Value = 0
$ws.Cells.Item(171, 17).Value = 0
$ws.Cells.Item(171, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 172
$ws.Cells.Item(172, 1).Value = 45495
$ws.Cells.Item(172, 2).Value = 1407
$ws.Cells.Item(172, 3).Value = 1464.050048828125
$ws.Cells.Item(172, 4).Value = 1361
$ws.Cells.Item(172, 5).Value = 1403
$ws.Cells.Item(172, 6).Value = 1400.541625976562
$ws.Cells.Item(172, 7).Value = 3737779
$ws.Cells.Item(172, 8).Value = 2024
$ws.Cells.Item(172, 9).Value = 7
$ws.Cells.Item(172, 10).Value = 22
$ws.Cells.Item(172, 11).Value = 0
$ws.Cells.Item(172, 12).Value = 0
$ws.Cells.Item(172, 13).Value = 0
$ws.Cells.Item(172, 14).Value = 30
$ws.Cells.Item(172, 15).Value = 0
$ws.Cells.Item(172, 16).Value = 0
$ws.Cells.Item(172, 17).Value = 0
$ws.Cells.Item(172, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 173
$ws.Cells.Item(173, 1).Value = 45502
$ws.Cells.Item(173, 2).Value = 1414.800048828125
$ws.Cells.Item(173, 3).Value = 1414.800048828125
$ws.Cells.Item(173, 4).Value = 1217.5
$ws.Cells.Item(173, 5).Value = 1231.150024414062
$ws.Cells.Item(173, 6).Value = 1228.99267578125
$ws.Cells.Item(173, 7).Value = 8059818
$ws.Cells.Item(173, 8).Value = 2024
$ws.Cells.Item(173, 9).Value = 7
$ws.Cells.Item(173, 10).Value = 29
$ws.Cells.Item(173, 11).Value = 0
$ws.Cells.Item(173, 12).Value = 0
$ws.Cells.Item(173, 13).Value = 0
$ws.Cells.Item(173, 14).Value = 31
$ws.Cells.Item(173, 15).Value = 0
$ws.Cells.Item(173, 16).Value = 0
$ws.Cells.Item(173, 17).Value = 2
$ws.Cells.Item(173, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 174
$ws.Cells.Item(174, 1).Value = 45509
$ws.Cells.Item(174, 2).Value = 1201.75
$ws.Cells.Item(174, 3).Value = 1275
$ws.Cells.Item(174, 4).Value = 1147.900024414062
$ws.Cells.Item(174, 5).Value = 1244.349975585938
$ws.Cells.Item(174, 6).Value = 1242.169555664062
$ws.Cells.Item(174, 7).Value = 8375809
$ws.Cells.Item(174, 8).Value = 2024
$ws.Cells.Item(174, 9).Value = 8
$ws.Cells.Item(174, 10).Value = 5
$ws.Cells.Item(174, 11).Value = 0
$ws.Cells.Item(174, 12).Value = 0
$ws.Cells.Item(174, 13).Value = 0
$ws.Cells.Item(174, 14).Value = 32
$ws.Cells.Item(174, 15).Value = 0
$ws.Cells.Item(174, 16).Value = 0
$ws.Cells.Item(174, 17).Value = 0
$ws.Cells.Item(174, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 175
$ws.Cells.Item(175, 1).Value = 45516
$ws.Cells.Item(175, 2).Value = 1242.25
$ws.Cells.Item(175, 3).Value = 1336
$ws.Cells.Item(175, 4).Value = 1235
$ws.Cells.Item(175, 5).Value = 1296.199951171875
$ws.Cells.Item(175, 6).Value = 1293.9287109375
$ws.Cells.Item(175, 7).Value = 5589373
$ws.Cells.Item(175, 8).Value = 2024
$ws.Cells.Item(175, 9).Value = 8
$ws.Cells.Item(175, 10).Value = 12
$ws.Cells.Item(175, 11).Value = 0
$ws.Cells.Item(175, 12).Value = 0
$ws.Cells.Item(175, 13).Value = 0
$ws.Cells.Item(175, 14).Value = 33
$ws.Cells.Item(175, 15).Value = 0
$ws.Cells.Item(175, 16).Value = 0
$ws.Cells.Item(175, 17).Value = 0
$ws.Cells.Item(175, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 176
$ws.Cells.Item(176, 1).Value = 45523
$ws.Cells.Item(176, 2).Value = 1297.400024414062
$ws.Cells.Item(176, 3).Value = 1327.5
$ws.Cells.Item(176, 4).Value = 1156.199951171875
$ws.Cells.Item(176, 5).Value = 1161.25
$ws.Cells.Item(176, 6).Value = 1161.25
$ws.Cells.Item(176, 7).Value = 3845638
$ws.Cells.Item(176, 8).Value = 2024
$ws.Cells.Item(176, 9).Value = 8
$ws.Cells.Item(176, 10).Value = 19
$ws.Cells.Item(176, 11).Value = 0
$ws.Cells.Item(176, 12).Value = 0
$ws.Cells.Item(176, 13).Value = 0
$ws.Cells.Item(176, 14).Value = 34
$ws.Cells.Item(176, 15).Value = 0
$ws.Cells.Item(176, 16).Value = 0
$ws.Cells.Item(176, 17).Value = 0
$ws.Cells.Item(176, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 177
$ws.Cells.Item(177, 1).Value = 45530
$ws.Cells.Item(177, 2).Value = 1161.75
$ws.Cells.Item(177, 3).Value = 1288.550048828125
$ws.Cells.Item(177, 4).Value = 1144
$ws.Cells.Item(177, 5).Value = 1253.900024414062
$ws.Cells.Item(177, 6).Value = 1253.900024414062
$ws.Cells.Item(177, 7).Value = 9462534
$ws.Cells.Item(177, 8).Value = 2024
$ws.Cells.Item(177, 9).Value = 8
$ws.Cells.Item(177, 10).Value = 26
$ws.Cells.Item(177, 11).Value = 0
$ws.Cells.Item(177, 12).Value = 0
$ws.Cells.Item(177, 13).Value = 0
$ws.Cells.Item(177, 14).Value = 35
$ws.Cells.Item(177, 15).Value = 0
$ws.Cells.Item(177, 16).Value = 0
$ws.Cells.Item(177, 17).Value = 0
$ws.Cells.Item(177, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 178
$ws.Cells.Item(178, 1).Value = 45537
$ws.Cells.Item(178, 2).Value = 1253.900024414062
$ws.Cells.Item(178, 3).Value = 1269
$ws.Cells.Item(178, 4).Value = 1181.849975585938
$ws.Cells.Item(178, 5).Value = 1186.849975585938
$ws.Cells.Item(178, 6).Value = 1186.849975585938
$ws.Cells.Item(178, 7).Value = 2945007
$ws.Cells.Item(178, 8).Value = 2024
$ws.Cells.Item(178, 9).Value = 9
$ws.Cells.Item(178, 10).Value = 2
$ws.Cells.Item(178, 11).Value = 0
$ws.Cells.Item(178, 12).Value = 0
$ws.Cells.Item(178, 13).Value = 0
$ws.Cells.Item(178, 14).Value = 36
$ws.Cells.Item(178, 15).Value = 0
$ws.Cells.Item(178, 16).Value = 0
$ws.Cells.Item(178, 17).Value = 0
$ws.Cells.Item(178, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 179
$ws.Cells.Item(179, 1).Value = 45544
$ws.Cells.Item(179, 2).Value = 1190.25
$ws.Cells.Item(179, 3).Value = 1247.699951171875
$ws.Cells.Item(179, 4).Value = 1165.599975585938
$ws.Cells.Item(179, 5).Value = 1230.199951171875
$ws.Cells.Item(179, 6).Value = 1230.199951171875
$ws.Cells.Item(179, 7).Value = 3480583
$ws.Cells.Item(179, 8).Value = 2024
$ws.Cells.Item(179, 9).Value = 9
$ws.Cells.Item(179, 10).Value = 9
$ws.Cells.Item(179, 11).Value = 0
$ws.Cells.Item(179, 12).Value = 0
$ws.Cells.Item(179, 13).Value = 0
$ws.Cells.Item(179, 14).Value = 37
$ws.Cells.Item(179, 15).Value = 0
$ws.Cells.Item(179, 16).Value = 0
$ws.Cells.Item(179, 17).Value = 0
$ws.Cells.Item(179, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 180
$ws.Cells.Item(180, 1).Value = 45551
$ws.Cells.Item(180, 2).Value = 1275
$ws.Cells.Item(180, 3).Value = 1410
$ws.Cells.Item(180, 4).Value = 1238.349975585938
$ws.Cells.Item(180, 5).Value = 1395.849975585938
$ws.Cells.Item(180, 6).Value = 1395.849975585938
$ws.Cells.Item(180, 7).Value = 8235705
$ws.Cells.Item(180, 8).Value = 2024
$ws.Cells.Item(180, 9).Value = 9
$ws.Cells.Item(180, 10).Value = 16
$ws.Cells.Item(180, 11).Value = 0
$ws.Cells.Item(180, 12).Value = 0
$ws.Cells.Item(180, 13).Value = 0
$ws.Cells.Item(180, 14).Value = 38
$ws.Cells.Item(180, 15).Value = 0
$ws.Cells.Item(180, 16).Value = 0
$ws.Cells.Item(180, 17).Value = 0
$ws.Cells.Item(180, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 181
$ws.Cells.Item(181, 1).Value = 45558
$ws.Cells.Item(181, 2).Value = 1419
$ws.Cells.Item(181, 3).Value = 1440
$ws.Cells.Item(181, 4).Value = 1280.849975585938
$ws.Cells.Item(181, 5).Value = 1300.400024414062
$ws.Cells.Item(181, 6).Value = 1300.400024414062
$ws.Cells.Item(181, 7).Value = 7786657
$ws.Cells.Item(181, 8).Value = 2024
$ws.Cells.Item(181, 9).Value = 9
$ws.Cells.Item(181, 10).Value = 23
$ws.Cells.Item(181, 11).Value = 0
$ws.Cells.Item(181, 12).Value = 0
$ws.Cells.Item(181, 13).Value = 0
$ws.Cells.Item(181, 14).Value = 39
$ws.Cells.Item(181, 15).Value = 1
$ws.Cells.Item(181, 16).Value = 0
$ws.Cells.Item(181, 17).Value = 0
$ws.Cells.Item(181, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 182
$ws.Cells.Item(182, 1).Value = 45565
$ws.Cells.Item(182, 2).Value = 1280
$ws.Cells.Item(182, 3).Value = 1288.949951171875
$ws.Cells.Item(182, 4).Value = 1123
$ws.Cells.Item(182, 5).Value = 1173.25
$ws.Cells.Item(182, 6).Value = 1173.25
$ws.Cells.Item(182, 7).Value = 8084580
$ws.Cells.Item(182, 8).Value = 2024
$ws.Cells.Item(182, 9).Value = 9
$ws.Cells.Item(182, 10).Value = 30
$ws.Cells.Item(182, 11).Value = 0
$ws.Cells.Item(182, 12).Value = 0
$ws.Cells.Item(182, 13).Value = 0
$ws.Cells.Item(182, 14).Value = 40
$ws.Cells.Item(182, 15).Value = 0
$ws.Cells.Item(182, 16).Value = 0
$ws.Cells.Item(182, 17).Value = 0
$ws.Cells.Item(182, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 183
$ws.Cells.Item(183, 1).Value = 45572
$ws.Cells.Item(183, 2).Value = 1204.949951171875
$ws.Cells.Item(183, 3).Value = 1225
$ws.Cells.Item(183, 4).Value = 1147.75
$ws.Cells.Item(183, 5).Value = 1169.849975585938
$ws.Cells.Item(183, 6).Value = 1169.849975585938
$ws.Cells.Item(183, 7).Value = 6910877
$ws.Cells.Item(183, 8).Value = 2024
$ws.Cells.Item(183, 9).Value = 10
$ws.Cells.Item(183, 10).Value = 7
$ws.Cells.Item(183, 11).Value = 0
$ws.Cells.Item(183, 12).Value = 0
$ws.Cells.Item(183, 13).Value = 0
$ws.Cells.Item(183, 14).Value = 41
$ws.Cells.Item(183, 15).Value = 0
$ws.Cells.Item(183, 16).Value = 0
$ws.Cells.Item(183, 17).Value = 0
$ws.Cells.Item(183, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 184
$ws.Cells.Item(184, 1).Value = 45579
$ws.Cells.Item(184, 2).Value = 1179
$ws.Cells.Item(184, 3).Value = 1234
$ws.Cells.Item(184, 4).Value = 1141.099975585938
$ws.Cells.Item(184, 5).Value = 1161
$ws.Cells.Item(184, 6).Value = 1161
$ws.Cells.Item(184, 7).Value = 4268357
$ws.Cells.Item(184, 8).Value = 2024
$ws.Cells.Item(184, 9).Value = 10
$ws.Cells.Item(184, 10).Value = 14
$ws.Cells.Item(184, 11).Value = 0
$ws.Cells.Item(184, 12).Value = 0
$ws.Cells.Item(184, 13).Value = 0
$ws.Cells.Item(184, 14).Value = 42
$ws.Cells.Item(184, 15).Value = 0
$ws.Cells.Item(184, 16).Value = 0
$ws.Cells.Item(184, 17).Value = 0
$ws.Cells.Item(184, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 185
$ws.Cells.Item(185, 1).Value = 45586
$ws.Cells.Item(185, 2).Value = 1167.849975585938
$ws.Cells.Item(185, 3).Value = 1181.949951171875
$ws.Cells.Item(185, 4).Value = 1043.5
$ws.Cells.Item(185, 5).Value = 1082.349975585938
$ws.Cells.Item(185, 6).Value = 1082.349975585938
$ws.Cells.Item(185, 7).Value = 5017519
$ws.Cells.Item(185, 8).Value = 2024
$ws.Cells.Item(185, 9).Value = 10
$ws.Cells.Item(185, 10).Value = 21
$ws.Cells.Item(185, 11).Value = 0
$ws.Cells.Item(185, 12).Value = 0
$ws.Cells.Item(185, 13).Value = 0
$ws.Cells.Item(185, 14).Value = 43
$ws.Cells.Item(185, 15).Value = 0
$ws.Cells.Item(185, 16).Value = 0
$ws.Cells.Item(185, 17).Value = 0
$ws.Cells.Item(185, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 186
$ws.Cells.Item(186, 1).Value = 45593
$ws.Cells.Item(186, 2).Value = 1111
$ws.Cells.Item(186, 3).Value = 1222
$ws.Cells.Item(186, 4).Value = 1043.050048828125
$ws.Cells.Item(186, 5).Value = 1207.449951171875
$ws.Cells.Item(186, 6).Value = 1207.449951171875
$ws.Cells.Item(186, 7).Value = 7259268
$ws.Cells.Item(186, 8).Value = 2024
$ws.Cells.Item(186, 9).Value = 10
$ws.Cells.Item(186, 10).Value = 28
$ws.Cells.Item(186, 11).Value = 0
$ws.Cells.Item(186, 12).Value = 0
$ws.Cells.Item(186, 13).Value = 0
$ws.Cells.Item(186, 14).Value = 44
$ws.Cells.Item(186, 15).Value = 2
$ws.Cells.Item(186, 16).Value = 0
$ws.Cells.Item(186, 17).Value = 0
$ws.Cells.Item(186, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 187
$ws.Cells.Item(187, 1).Value = 45600
$ws.Cells.Item(187, 2).Value = 1212.75
$ws.Cells.Item(187, 3).Value = 1285.449951171875
$ws.Cells.Item(187, 4).Value = 1152.050048828125
$ws.Cells.Item(187, 5).Value = 1170.300048828125
$ws.Cells.Item(187, 6).Value = 1170.300048828125
$ws.Cells.Item(187, 7).Value = 5508708
$ws.Cells.Item(187, 8).Value = 2024
$ws.Cells.Item(187, 9).Value = 11
$ws.Cells.Item(187, 10).Value = 4
$ws.Cells.Item(187, 11).Value = 0
$ws.Cells.Item(187, 12).Value = 0
$ws.Cells.Item(187, 13).Value = 0
$ws.Cells.Item(187, 14).Value = 45
$ws.Cells.Item(187, 15).Value = 0
$ws.Cells.Item(187, 16).Value = 0
$ws.Cells.Item(187, 17).Value = 2
$ws.Cells.Item(187, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 188
$ws.Cells.Item(188, 1).Value = 45607
$ws.Cells.Item(188, 2).Value = 1170.25
$ws.Cells.Item(188, 3).Value = 1291
$ws.Cells.Item(188, 4).Value = 1161.75
$ws.Cells.Item(188, 5).Value = 1222.150024414062
$ws.Cells.Item(188, 6).Value = 1222.150024414062
$ws.Cells.Item(188, 7).Value = 7519820
$ws.Cells.Item(188, 8).Value = 2024
$ws.Cells.Item(188, 9).Value = 11
$ws.Cells.Item(188, 10).Value = 11
$ws.Cells.Item(188, 11).Value = 0
$ws.Cells.Item(188, 12).Value = 0
$ws.Cells.Item(188, 13).Value = 0
$ws.Cells.Item(188, 14).Value = 46
$ws.Cells.Item(188, 15).Value = 0
$ws.Cells.Item(188, 16).Value = 0
$ws.Cells.Item(188, 17).Value = 0
$ws.Cells.Item(188, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 189
$ws.Cells.Item(189, 1).Value = 45614
$ws.Cells.Item(189, 2).Value = 1263
$ws.Cells.Item(189, 3).Value = 1294
$ws.Cells.Item(189, 4).Value = 1163.050048828125
$ws.Cells.Item(189, 5).Value = 1244.849975585938
$ws.Cells.Item(189, 6).Value = 1244.849975585938
$ws.Cells.Item(189, 7).Value = 7414855
$ws.Cells.Item(189, 8).Value = 2024
$ws.Cells.Item(189, 9).Value = 11
$ws.Cells.Item(189, 10).Value = 18
$ws.Cells.Item(189, 11).Value = 0
$ws.Cells.Item(189, 12).Value = 0
$ws.Cells.Item(189, 13).Value = 0
$ws.Cells.Item(189, 14).Value = 47
$ws.Cells.Item(189, 15).Value = 0
$ws.Cells.Item(189, 16).Value = 0
$ws.Cells.Item(189, 17).Value = 0
$ws.Cells.Item(189, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 190
$ws.Cells.Item(190, 1).Value = 45621
$ws.Cells.Item(190, 2).Value = 1324.949951171875
$ws.Cells.Item(190, 3).Value = 1326.400024414062
$ws.Cells.Item(190, 4).Value = 1212.150024414062
$ws.Cells.Item(190, 5).Value = 1253.199951171875
$ws.Cells.Item(190, 6).Value = 1253.199951171875
$ws.Cells.Item(190, 7).Value = 12007210
$ws.Cells.Item(190, 8).Value = 2024
$ws.Cells.Item(190, 9).Value = 11
$ws.Cells.Item(190, 10).Value = 25
$ws.Cells.Item(190, 11).Value = 0
$ws.Cells.Item(190, 12).Value = 0
$ws.Cells.Item(190, 13).Value = 0
$ws.Cells.Item(190, 14).Value = 48
$ws.Cells.Item(190, 15).Value = 0
$ws.Cells.Item(190, 16).Value = 0
$ws.Cells.Item(190, 17).Value = 0
$ws.Cells.Item(190, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 191
$ws.Cells.Item(191, 1).Value = 45628
$ws.Cells.Item(191, 2).Value = 1254.449951171875
$ws.Cells.Item(191, 3).Value = 1397.400024414062
$ws.Cells.Item(191, 4).Value = 1246
$ws.Cells.Item(191, 5).Value = 1371.300048828125
$ws.Cells.Item(191, 6).Value = 1371.300048828125
$ws.Cells.Item(191, 7).Value = 9695039
$ws.Cells.Item(191, 8).Value = 2024
$ws.Cells.Item(191, 9).Value = 12
$ws.Cells.Item(191, 10).Value = 2
$ws.Cells.Item(191, 11).Value = 0
$ws.Cells.Item(191, 12).Value = 0
$ws.Cells.Item(191, 13).Value = 0
$ws.Cells.Item(191, 14).Value = 49
$ws.Cells.Item(191, 15).Value = 0
$ws.Cells.Item(191, 16).Value = 0
$ws.Cells.Item(191, 17).Value = 0
$ws.Cells.Item(191, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
